$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.543.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.14%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.804.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.88%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.56%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  +0.57%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'308.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.03%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4540"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3668"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.51%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07140"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.51%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8704"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.18%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07779"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.67%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -3.79%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.808.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.31%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -2.38%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.325"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.79%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'86.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.83%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +0.66%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008584"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.54%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.009"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'26.583.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.67%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -3.47%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.068.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.63%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.98%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.985"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.41%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'151.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.81%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'17.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.81%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.987"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.57%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'113.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.75%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -4.56%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08699"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.76%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.030"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.88%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.7335"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.60%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.437"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.75%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -5.55%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.488"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -7.54%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.079"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.36%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01918"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.32%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05093"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.98%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.859"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.61%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.871"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.04%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4905"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.90%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -4.49%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.103"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.64%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.66%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4591"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.77%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'102.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.48%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'9.964"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.47%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.583"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.38%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.68%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'63.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.92%  "
$ws.Range("E51").Style = "Normal"
